$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (C1:F1 are brand new, G1 takes over the old "ExpPoints" header) ---
# Copy the formatting of the existing bold/bordered header style (A1) onto the
# new header cells before writing their text, so they pick up the same style
# index instead of minting new (duplicate) styles.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1:G1").PasteSpecial(-4122) | Out-Null

$ws.Range("G1").Value = "ExpPoints"
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "RELEGATION"

# --- Re-ranked table body: Rank (A) / Team (B) / ExpPoints (now G) ---
# Manchester United and Tottenham Hotspur swap (rows 10/12), and
# Nottingham Forest / Burnley swap (rows 18/19), per the refreshed model.
$data = @(
    @(1,  'Arsenal',                  81.56676736176716),
    @(2,  'Manchester City',          71.68842760627921),
    @(3,  'Liverpool',                68.49771028076655),
    @(4,  'Chelsea',                  60.54405966664068),
    @(5,  'Aston Villa',              59.8921389461057),
    @(6,  'Crystal Palace',           57.99553001686682),
    @(7,  'AFC Bournemouth',          56.84799701114331),
    @(8,  'Newcastle United',         55.44259210585084),
    @(9,  'Manchester United',        54.45524036993888),
    @(10, 'Brighton & Hove Albion',   54.4523841883292),
    @(11, 'Tottenham Hotspur',        54.06132578230342),
    @(12, 'Brentford',                52.49121651768586),
    @(13, 'Sunderland',               44.12682558533493),
    @(14, 'Everton',                  44.01075920617485),
    @(15, 'Fulham',                   42.15728998551921),
    @(16, 'Leeds United',             37.43795281308542),
    @(17, 'Nottingham Forest',        36.45785696161575),
    @(18, 'Burnley',                  35.9863376514443),
    @(19, 'West Ham United',          33.38041117972701),
    @(20, 'Wolverhampton Wanderers',  28.35410445299532)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rank = $data[$i][0]
    $team = $data[$i][1]
    $pts  = $data[$i][2]

    $ws.Cells.Item($r, 1).Value = $rank
    $ws.Cells.Item($r, 2).Value = $team

    # New placeholder columns (WIN / TOP4 / TOP5 / RELEGATION) - empty text
    # cells reserved for the upcoming Monte-Carlo simulation outputs.
    $ws.Cells.Item($r, 3).Value = "'"
    $ws.Cells.Item($r, 4).Value = "'"
    $ws.Cells.Item($r, 5).Value = "'"
    $ws.Cells.Item($r, 6).Value = "'"
    $ws.Cells.Item($r, 3).Style = "Normal"
    $ws.Cells.Item($r, 4).Style = "Normal"
    $ws.Cells.Item($r, 5).Style = "Normal"
    $ws.Cells.Item($r, 6).Style = "Normal"

    # ExpPoints now lives in column G
    $ws.Cells.Item($r, 7).Value = $pts
}
